# Update projects dan publications
# - Row 2 becomes the "Improving Quality of Disease Preparedness... (INSPIRASI)" project,
#   with period "2021-present"
# - Row 3 becomes the "Integrated Health Surveillance System... / World Bank Group" project,
#   with period "2021-2022"
# (previously these two projects occupied rows 3 and 2 respectively, with a plain
# numeric year of 2021 in column A for both)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$b2 = "Improving Quality of Disease Preparedness, Surveillance & Response In Indonesia (INSPIRASI) / CDC USA"
$c2 = "Researcher : Analyse Indonesia national health insurance claims on capability to monitor potentially outbreak infectious disease surveillance using time series analysis and data visualization. Designing dashboard presentation of infectious disease surveillance using Indonesia national health insurance claims data to support potentially outbreak infectious disease surveillance program in Indonesia. Writing scientific publication on quality of Early Warning and Response System reporting in primary health care facility. Writing policy brief on implementation of Public Health Emergency Operation Center in regional settings."

$b3 = "Integrated Health Surveillance System Supporting Covid-19 Emergency Response Operations / World Bank Group"
$c3 = "Consultant: Analyse workflow of Allrecord TC-19, an application used by Indonesia government to collect PCR and rapid antigen result for COVID-19 surveillance. Building online user guideline for Allrecord TC-19 application to help user get the update from the app. Held two batch of one week FHIR workshops for public/private hospitals, HIMS vendors, and healthcare startups."

# Column A: switch from plain year numbers to explicit period ranges.
# Intern "2021-2022" before "2021-present" so the shared-string table
# order matches the source workbook.
$ws.Range("A3").Value = "2021-2022"
$ws.Range("A2").Value = "2021-present"

# Column B / C: swap so row 2 holds the INSPIRASI project and row 3 holds
# the World Bank Group project
$ws.Range("B2").Value = $b2
$ws.Range("C2").Value = $c2

$ws.Range("B3").Value = $b3
$ws.Range("C3").Value = $c3

# Restore the selection to B2 (matches the author's saved view state)
$ws.Range("B2").Select()
